# Update training metrics on the "con_opt_hip" worksheet to reflect the
# latest training run results (Se actualiza Readme y artefactos).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("con_opt_hip")

# epoch 1
$ws.Range("B2").Value = 0.1543
$ws.Range("C2").Value = 0.0146
$ws.Range("D2").Value = 0.9978
$ws.Range("E2").Value = 0.9978

# epoch 2
$ws.Range("B3").Value = 0.0147
$ws.Range("C3").Value = 0.0203
$ws.Range("D3").Value = 0.996
$ws.Range("E3").Value = 0.996

# epoch 3
$ws.Range("B4").Value = 0.0102
$ws.Range("C4").Value = 0.0144
$ws.Range("D4").Value = 0.9971
$ws.Range("E4").Value = 0.9971

# epoch 4
$ws.Range("B5").Value = 0.124
$ws.Range("C5").Value = 0.0671
$ws.Range("D5").Value = 0.9914
$ws.Range("E5").Value = 0.9914

# epoch 5
$ws.Range("B6").Value = 0.0205
$ws.Range("C6").Value = 0.0219
$ws.Range("D6").Value = 0.996
$ws.Range("E6").Value = 0.996

# epoch 6
$ws.Range("B7").Value = 0.0083
$ws.Range("C7").Value = 0.0188
$ws.Range("D7").Value = 0.9964
$ws.Range("E7").Value = 0.9964

# epoch 7
$ws.Range("B8").Value = 0.0039
$ws.Range("C8").Value = 0.24
$ws.Range("D8").Value = 0.996
$ws.Range("E8").Value = 0.996

# epoch 8
$ws.Range("B9").Value = 0.0109
$ws.Range("C9").Value = 0.0194
$ws.Range("D9").Value = 0.9952
$ws.Range("E9").Value = 0.9952

# epoch 9
$ws.Range("B10").Value = 0.0098
$ws.Range("C10").Value = 0.0246
$ws.Range("D10").Value = 0.9963
$ws.Range("E10").Value = 0.9963

# epoch 10
$ws.Range("B11").Value = 0.0061
$ws.Range("C11").Value = 0.0206
$ws.Range("D11").Value = 0.9963
$ws.Range("E11").Value = 0.9963
